$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New rows 18-21: clone the formatting of row 17 (same border/number-format
#    pattern used by every data row), then fix up H20:H21 which use the plain
#    "General" style instead of the 0.0 speed-up format.
# ---------------------------------------------------------------------------
$ws.Range("A17:K17").Copy()
$ws.Range("A18:K19").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A20:K21").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A2").Copy()
$ws.Range("H20:H21").PasteSpecial(-4122)   # xlPasteFormats (plain style, like column A)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) New row values (A:I) - practice 2, part 3: ICC compiler variant of
#    axpy_stride, plus the new "cond" benchmark (cond_esc / cond_vec).
# ---------------------------------------------------------------------------
$ws.Range("A18").Value = "axpy_stride"
$ws.Range("B18").Value = "axpy_stride_v1()"
$ws.Range("C18").Value = "ICC"
$ws.Range("D18").Value = 2000
$ws.Range("E18").Value = "Single"
$ws.Range("F18").Value = 3.54
$ws.Range("G18").Value = "-"
$ws.Range("H18").Value = 1.8
$ws.Range("I18").Value = 478502976

$ws.Range("A19").Value = "axpy_stride"
$ws.Range("B19").Value = "axpy_stride_v2(x, y)"
$ws.Range("C19").Value = "ICC"
$ws.Range("D19").Value = 2000
$ws.Range("E19").Value = "Single"
$ws.Range("F19").Value = 3.54
$ws.Range("G19").Value = "-"
$ws.Range("H19").Value = 1.8
$ws.Range("I19").Value = 478502976

$ws.Range("A20").Value = "cond"
$ws.Range("B20").Value = "cond_esc"
$ws.Range("C20").Value = "GCC"
$ws.Range("D20").Value = 2000
$ws.Range("E20").Value = "Single"
$ws.Range("F20").Value = 5.5
$ws.Range("G20").Value = "-"
$ws.Range("H20").Value = 1.1
$ws.Range("I20").Value = 102.3685

$ws.Range("A21").Value = "cond"
$ws.Range("B21").Value = "cond_vec"
$ws.Range("C21").Value = "GCC"
$ws.Range("D21").Value = 2000
$ws.Range("E21").Value = "Single"
$ws.Range("F21").Value = 1.6
$ws.Range("G21").Value = "-"
$ws.Range("H21").Value = 3.9
$ws.Range("I21").Value = 1023.693359

# ---------------------------------------------------------------------------
# 3) Column J (Speedup vs the reference variant): every group of related
#    variants is re-based against its own reference implementation.
#    Rewriting re-creates the "shared formula" groups with the right
#    boundaries (Excel groups a formula typed once across a multi-cell range).
# ---------------------------------------------------------------------------
$ws.Range("J2").Formula = "=`$F`$5/F2"
$ws.Range("J3:J11").Formula = "=`$F`$5/F3"
$ws.Range("J12").Formula = "=`$F`$12/F12"
$ws.Range("J13:J15").Formula = "=`$F`$12/F13"
$ws.Range("J16").Formula = "=`$F`$16/F16"
$ws.Range("J17").Formula = "=`$F`$16/F17"
$ws.Range("J18:J19").Formula = "=`$F`$16/F18"
$ws.Range("J20").Formula = "=`$F`$20/F20"
$ws.Range("J21").Formula = "=`$F`$20/F21"

# ---------------------------------------------------------------------------
# 4) Column K (GFLOPS) for the new rows follows the same pattern already used
#    for rows 3-17.
# ---------------------------------------------------------------------------
$ws.Range("K18:K21").Formula = "=D18/100/F18"

# ---------------------------------------------------------------------------
# 5) Selection, to match where the author left the cursor.
# ---------------------------------------------------------------------------
$ws.Range("E19").Select()
